$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The annotator duplicated row 55 into a brand-new row 56 (same Annotator /
# politeness_score / polite_expressions / issue_type), then edited the
# sentence_purpose, id, source_file and text of the new row. Copy/paste the
# existing row first so the duplicated cells keep row 55's original text
# formatting (politeness_score "3" stored as text), then tweak B55 to the
# numeric 3 it ends up with, and fill in row 56's unique fields.
$ws.Range("A55:H55").Copy()
$ws.Range("A56:H56").PasteSpecial()

$ws.Range("B55").Value = 3

$ws.Range("D56").Value = "SMY"
$ws.Range("F56").Value = "e3eeb88f-0832-4aa9-a6cc-39ada0451b32"
$ws.Range("G56").Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Range("H56").Value = "This paper shows that models trained on a synthetic dataset are vulnerable to small adversarial perturbations which lie on the data manifold."
